$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "yyyy\-mm\-dd;@"

# Data for the new "fall 2024" term block (rows 50-64), mirroring the
# existing "fall 2023" block's layout: week_start, week_number, date,
# exam/holiday label (D), exam/holiday date (E), holiday name (G), topic (H)
$rows = @(
    @{ R=50; WK=1;  Date=45530; Topic="Python Basics" }
    @{ R=51; WK=2;  Date=45539; Topic="Functions, Operators and Expressions" }
    @{ R=52; WK=3;  Date=45546; EDate=45171; EBlackFont=$true; Holiday="Labor Day"; Topic="Functions, Operators and Expressions" }
    @{ R=53; WK=4;  Date=45553; Topic="Control Flow (if statements)" }
    @{ R=54; WK=5;  Date=45560; Label="Midterm 1"; EDate=45194; Topic="Control Flow (while)" }
    @{ R=55; WK=6;  Date=45569; Topic="Data Structures (lists)" }
    @{ R=56; WK=7;  Date=45576; Topic="Control Flow (for loops)" }
    @{ R=57; WK=8;  Date=45583; Topic="Control Flow (for loops)" }
    @{ R=58; WK=9;  Date=45590; Label="Midterm 2"; EDate=45222; Topic="Data Structures (dictionaries)" }
    @{ R=59; WK=10; Date=45597; Topic="Input and Output" }
    @{ R=60; WK=11; Date=45604; EDate=45240; Holiday="Veterans Day"; Topic="Data Structures (tuples)" }
    @{ R=61; WK=12; Date=45609; Topic="Data Structures (sets)" }
    @{ R=62; WK=13; Date=45616; Label="Midterm 3"; EDate=45616; Topic="Mutability" }
    @{ R=63; WK=14; Date=45628; EDate=45258; Holiday="Thanksgiving Recess"; Topic="Control Flow + Data Structures" }
    @{ R=64; WK=15; Date=45635; Label="Final Exam"; EText="TBA"; Time="6:00pm - 8:00pm"; Topic="Review" }
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Cells.Item($r, 1).Value = "fall 2024"
    $ws.Cells.Item($r, 2).Value = $row.WK

    $ws.Cells.Item($r, 3).Value = $row.Date
    $ws.Cells.Item($r, 3).NumberFormat = $dateFmt

    if ($row.Label) {
        $ws.Cells.Item($r, 4).Value = $row.Label
    }

    if ($row.EDate) {
        $ws.Cells.Item($r, 5).Value = $row.EDate
        $ws.Cells.Item($r, 5).NumberFormat = $dateFmt
        if ($row.EBlackFont) {
            $ws.Cells.Item($r, 5).Font.Color = 0
        }
    } elseif ($row.EText) {
        $ws.Cells.Item($r, 5).Value = $row.EText
        $ws.Cells.Item($r, 5).NumberFormat = $dateFmt
    }

    if ($row.Time) {
        $ws.Cells.Item($r, 6).Value = $row.Time
    }

    if ($row.Holiday) {
        $ws.Cells.Item($r, 7).Value = $row.Holiday
    }

    $ws.Cells.Item($r, 8).Value = $row.Topic
}

$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D41").Select()
